# Fruta / hortaliza, semanal
# Adds a new weekly price record for "Vega Modelo de Temuco - Camote".
# The new record is inserted as row 33 (pushing the existing rows 33-73
# down to 34-74, growing the sheet from A1:R73 to A1:R74).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 33, shifting everything below
# (rows 33-73) down by one.
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with the new weekly record. Values
# that are identical across the whole "Vega Modelo de Temuco - Camote"
# block (Mercado ID, Mercado, Región, Codreg, Categoría ID, Categoría,
# Variedad, Calidad, Precio mínimo/máximo/promedio, Unidad de
# comercialización, Origen, Precio $/Kg, Kg o Unidades, Clasificación)
# keep the same values as the row that used to be row 33; only the date
# (Fecha) and Volumen differ for this new entry.
$ws.Range("A33").Value = 10
$ws.Range("B33").Value = "Vega Modelo de Temuco"
$ws.Range("C33").Value = "La Araucanía"
$ws.Range("D33").Value = 44665
$ws.Range("E33").Value = 9
$ws.Range("F33").Value = 100114002
$ws.Range("G33").Value = "Camote"
$ws.Range("H33").Value = "Sin especificar"
$ws.Range("I33").Value = "Primera"
$ws.Range("J33").Value = 40
$ws.Range("K33").Value = 18000
$ws.Range("L33").Value = 18000
$ws.Range("M33").Value = 18000
$ws.Range("N33").Value = "$/malla 20 kilos"
$ws.Range("O33").Value = "Perú"
$ws.Range("P33").Value = 900
$ws.Range("Q33").Value = 20
$ws.Range("R33").Value = "Hortaliza"
